$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "E2"=3; "G2"=43.995596; "H2"=131.986788; "I2"=0.08241811124115486; "J2"=0.08241811124115485; "K2"=3; "M2"=7.214110666666667; "N2"=21.642332; "O2"=0.4688823795981188; "P2"=0.4688823795981188; "Q2"=317.3890983899573; "R2"=2856.501885509616; "S2"=0.03864440012073516; "T2"=0.03864440012073515
    "E3"=3; "G3"=43.995596; "H3"=131.986788; "I3"=0.08241811124115486; "J3"=0.08241811124115485; "K3"=3; "M3"=7.110350666666666; "N3"=21.331052; "O3"=0.4621384803214003; "P3"=0.4621384803214003; "Q3"=312.8241153489973; "R3"=2815.417038140976; "S3"=0.03808858067994743; "T3"=0.03808858067994742
    "E4"=3; "G4"=43.995596; "H4"=131.986788; "I4"=0.08241811124115486; "J4"=0.08241811124115485; "K4"=3; "M4"=1.061296333333333; "N4"=3.183889; "O4"=0.06897914008048092; "P4"=0.06897914008048092; "Q4"=46.69236471761466; "R4"=420.231282458532; "S4"=0.00568513044047228; "T4"=0.005685130440472279
    "E5"=3; "G5"=439.8208616666666; "H5"=1319.462585; "I5"=0.8239280291378236; "J5"=0.8239280291378236; "K5"=3; "M5"=7.214110666666667; "N5"=21.642332; "O5"=0.4688823795981188; "P5"=0.4688823795981188; "Q5"=3172.916369572024; "R5"=28556.24732614822; "S5"=0.3863253349197309; "T5"=0.3863253349197309
    "E6"=3; "G6"=439.8208616666666; "H6"=1319.462585; "I6"=0.8239280291378236; "J6"=0.8239280291378236; "K6"=3; "M6"=7.110350666666666; "N6"=21.331052; "O6"=0.4621384803214003; "P6"=0.4621384803214003; "Q6"=3127.280556965491; "R6"=28145.52501268942; "S6"=0.3807688472799602; "T6"=0.3807688472799602
    "E7"=3; "G7"=439.8208616666666; "H7"=1319.462585; "I7"=0.8239280291378236; "J7"=0.8239280291378236; "K7"=3; "M7"=1.061296333333333; "N7"=3.183889; "O7"=0.06897914008048092; "P7"=0.06897914008048092; "Q7"=466.7802678103405; "R7"=4201.022410293064; "S7"=0.0568338469381325; "T7"=0.0568338469381325
    "E8"=3; "G8"=49.99334866666667; "H8"=149.980046; "I8"=0.09365385962102149; "J8"=0.09365385962102149; "K8"=3; "M8"=7.214110666666667; "N8"=21.642332; "O8"=0.4688823795981188; "P8"=0.4688823795981188; "Q8"=360.6575498785858; "R8"=3245.917948907272; "S8"=0.04391264455765273; "T8"=0.04391264455765273
    "E9"=3; "G9"=49.99334866666667; "H9"=149.980046; "I9"=0.09365385962102149; "J9"=0.09365385962102149; "K9"=3; "M9"=7.110350666666666; "N9"=21.331052; "O9"=0.4621384803214003; "P9"=0.4621384803214003; "Q9"=355.4702400209325; "R9"=3199.232160188392; "S9"=0.04328105236149263; "T9"=0.04328105236149263
    "E10"=3; "G10"=49.99334866666667; "H10"=149.980046; "I10"=0.09365385962102149; "J10"=0.09365385962102149; "K10"=3; "M10"=1.061296333333333; "N10"=3.183889; "O10"=0.06897914008048092; "P10"=0.06897914008048092; "Q10"=53.05775763098822; "R10"=477.519818678894; "S10"=0.006460162701876137; "T10"=0.006460162701876137
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

Write-Output "Updated $($updates.Count) cells"